$d = $word.ActiveDocument

$d.Content.Find.Execute("386÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "691÷5=", 2) | Out-Null
$d.Content.Find.Execute("735÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "463÷3=", 2) | Out-Null
$d.Content.Find.Execute("681÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "226÷7=", 2) | Out-Null
$d.Content.Find.Execute("804÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "814÷6=", 2) | Out-Null
$d.Content.Find.Execute("998÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷3=", 2) | Out-Null
$d.Content.Find.Execute("190÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "760÷2=", 2) | Out-Null
$d.Content.Find.Execute("759÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "652÷2=", 2) | Out-Null
$d.Content.Find.Execute("116÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "778÷6=", 2) | Out-Null
$d.Content.Find.Execute("148÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "658÷6=", 2) | Out-Null
$d.Content.Find.Execute("869÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "951÷3=", 2) | Out-Null
$d.Content.Find.Execute("289÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "139÷4=", 2) | Out-Null
$d.Content.Find.Execute("847÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "256÷5=", 2) | Out-Null
$d.Content.Find.Execute("170÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "665÷6=", 2) | Out-Null
$d.Content.Find.Execute("102÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "475÷4=", 2) | Out-Null
$d.Content.Find.Execute("618÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "189÷6=", 2) | Out-Null
$d.Content.Find.Execute("338÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "770÷6=", 2) | Out-Null
$d.Content.Find.Execute("628÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "780÷6=", 2) | Out-Null
$d.Content.Find.Execute("568÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "735÷6=", 2) | Out-Null
$d.Content.Find.Execute("930÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "342÷7=", 2) | Out-Null
$d.Content.Find.Execute("954÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "226÷6=", 2) | Out-Null
$d.Content.Find.Execute("199÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "477÷2=", 2) | Out-Null
$d.Content.Find.Execute("767÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "314÷7=", 2) | Out-Null
$d.Content.Find.Execute("555÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "398÷7=", 2) | Out-Null
$d.Content.Find.Execute("759÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "476÷8=", 2) | Out-Null
$d.Content.Find.Execute("124÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "291÷6=", 2) | Out-Null
